$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-target the existing conditional formatting rules onto the columns they
# will occupy once C:E are removed (column deletion does not, by itself,
# shift the sqref of conditional formats in this engine), keeping their
# original dxfId/priority/dxf entries intact.
$quote = [char]34

$rngFailedSuccess = $ws.Range("F1:F1048576")
$fcFailed = $rngFailedSuccess.FormatConditions.Item(1)
$fcFailed.ModifyAppliesToRange($ws.Range("C1:C1048576"))
$ws.Range("C1:C1048576").FormatConditions.Item(1).Formula1 = "=NOT(ISERROR(SEARCH(" + $quote + "Failed" + $quote + ",C1)))"

$fcSuccessful = $ws.Range("C1:C1048576").FormatConditions.Item(2)
$fcSuccessful.ModifyAppliesToRange($ws.Range("C1:C1048576"))
$ws.Range("C1:C1048576").FormatConditions.Item(2).Formula1 = "=NOT(ISERROR(SEARCH(" + $quote + "Successful" + $quote + ",C1)))"

$rngFailureMessage = $ws.Range("H51:H1048576")
$fcFailureMessage = $rngFailureMessage.FormatConditions.Item(1)
$fcFailureMessage.ModifyAppliesToRange($ws.Range("E51:E1048576"))
$ws.Range("E51:E1048576").FormatConditions.Item(1).Formula1 = "=NOT(ISERROR(SEARCH(" + $quote + "Failure Message" + $quote + ",E51)))"

# Remove the "Report Type", "Download Name" and "New Name" columns (C:E).
# This shifts Status/Note/(blank) left from F:H into C:E, and also shifts
# all the following (now unused) column width definitions left by 3.
$ws.Range("C:E").Delete()

# Rename the first two remaining header columns.
$ws.Range("A1").Value = "FileName"
$ws.Range("B1").Value = "Date"
